$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels: two new "IMF (20%)" columns are inserted before the
# existing IMF columns, and the existing IMF data is pushed two columns to
# the right. The old "OECD (20%)" columns are removed entirely, so the
# trailing OECD columns stay where they were.
$ws.Range("F1").Value = "IMF (20%) - Sales"
$ws.Range("G1").Value = "IMF (20%) - Sales + Emp"
$ws.Range("H1").Value = "IMF - Sales"
$ws.Range("I1").Value = "IMF - Sales + Emp"

# Row 2 - shift old IMF values (F2,G2) into H2,I2, then set new IMF(20%) values
$f2 = $ws.Range("F2").Value2
$g2 = $ws.Range("G2").Value2
$ws.Range("H2").Value2 = $f2
$ws.Range("I2").Value2 = $g2
$ws.Range("F2").Value2 = 0.7444774198808549
$ws.Range("G2").Value2 = 0.606637845192245

# Row 3
$f3 = $ws.Range("F3").Value2
$g3 = $ws.Range("G3").Value2
$ws.Range("H3").Value2 = $f3
$ws.Range("I3").Value2 = $g3
$ws.Range("F3").Value2 = 2.973125562628548
$ws.Range("G3").Value2 = 4.34634358097193

# Row 4
$f4 = $ws.Range("F4").Value2
$g4 = $ws.Range("G4").Value2
$ws.Range("H4").Value2 = $f4
$ws.Range("I4").Value2 = $g4
$ws.Range("F4").Value2 = -0.05012849520288268
$ws.Range("G4").Value2 = 1.716044213773076

# Row 5
$f5 = $ws.Range("F5").Value2
$g5 = $ws.Range("G5").Value2
$ws.Range("H5").Value2 = $f5
$ws.Range("I5").Value2 = $g5
$ws.Range("F5").Value2 = -6.426616309521044
$ws.Range("G5").Value2 = -8.563732635657207

# Row 6
$f6 = $ws.Range("F6").Value2
$g6 = $ws.Range("G6").Value2
$ws.Range("H6").Value2 = $f6
$ws.Range("I6").Value2 = $g6
$ws.Range("F6").Value2 = 0.2446805784062994
$ws.Range("G6").Value2 = 0.5593890638537687
